# Apply updated cash high/low figures to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 2612.4499999999998
$ws.Range("C2").Value = 2531.1
$ws.Range("D2").Value = 2542.1999999999998
$ws.Range("E2").Value = 2539.5500000000002
$ws.Range("F2").Value = 45
$ws.Range("G2").Value = 2571.65
$ws.Range("B3").Value = 399.55
$ws.Range("C3").Value = 393.1
$ws.Range("D3").Value = 394.7
$ws.Range("E3").Value = 395.6
$ws.Range("F3").Value = 17
$ws.Range("G3").Value = 397.5
$ws.Range("B4").Value = 1492
$ws.Range("C4").Value = 1465
$ws.Range("D4").Value = 1467
$ws.Range("E4").Value = 1468.6
$ws.Range("F4").Value = 8
$ws.Range("G4").Value = 1489.2
$ws.Range("B5").Value = 7185.95
$ws.Range("C5").Value = 7121
$ws.Range("D5").Value = 7124
$ws.Range("E5").Value = 7137.45
$ws.Range("F5").Value = 9
$ws.Range("G5").Value = 7179
$ws.Range("B6").Value = 240.7
$ws.Range("C6").Value = 234
$ws.Range("D6").Value = 235.2
$ws.Range("E6").Value = 234.9
$ws.Range("F6").Value = 111
$ws.Range("G6").Value = 240
$ws.Range("B7").Value = 194.95
$ws.Range("C7").Value = 190.75
$ws.Range("D7").Value = 191.9
$ws.Range("E7").Value = 191.85
$ws.Range("F7").Value = 159
$ws.Range("G7").Value = 194.55
$ws.Range("B8").Value = 232.65
$ws.Range("C8").Value = 228.5
$ws.Range("D8").Value = 229.75
$ws.Range("E8").Value = 229.25
$ws.Range("F8").Value = 63
$ws.Range("G8").Value = 232.3
$ws.Range("B9").Value = 490.75
$ws.Range("C9").Value = 483
$ws.Range("D9").Value = 483.5
$ws.Range("E9").Value = 483.55
$ws.Range("F9").Value = 26
$ws.Range("G9").Value = 484.85
$ws.Range("B10").Value = 3367.7
$ws.Range("C10").Value = 3320
$ws.Range("D10").Value = 3325
$ws.Range("E10").Value = 3341.65
$ws.Range("F10").Value = 8
$ws.Range("G10").Value = 3353
$ws.Range("B11").Value = 143.15
$ws.Range("C11").Value = 139.4
$ws.Range("D11").Value = 140.6
$ws.Range("E11").Value = 140.69999999999999
$ws.Range("F11").Value = 226
$ws.Range("G11").Value = 142.80000000000001
$ws.Range("B12").Value = 1188.8
$ws.Range("C12").Value = 1169
$ws.Range("D12").Value = 1172.95
$ws.Range("E12").Value = 1171.25
$ws.Range("F12").Value = 24
$ws.Range("G12").Value = 1187
$ws.Range("B13").Value = 1596.2
$ws.Range("C13").Value = 1576.15
$ws.Range("D13").Value = 1579.6
$ws.Range("E13").Value = 1579.3
$ws.Range("F13").Value = 215
$ws.Range("G13").Value = 1591.05
$ws.Range("B14").Value = 464.65
$ws.Range("C14").Value = 453.85
$ws.Range("D14").Value = 455.55
$ws.Range("E14").Value = 456.5
$ws.Range("F14").Value = 59
$ws.Range("G14").Value = 464.4
$ws.Range("B15").Value = 982.4
$ws.Range("C15").Value = 966.75
$ws.Range("D15").Value = 968.8
$ws.Range("E15").Value = 968.95
$ws.Range("F15").Value = 283
$ws.Range("G15").Value = 974.3
$ws.Range("B16").Value = 1435
$ws.Range("C16").Value = 1411.95
$ws.Range("D16").Value = 1418.65
$ws.Range("E16").Value = 1425.15
$ws.Range("F16").Value = 35
$ws.Range("G16").Value = 1414.15
$ws.Range("B17").Value = 1425.05
$ws.Range("C17").Value = 1414.45
$ws.Range("D17").Value = 1424
$ws.Range("E17").Value = 1423.6
$ws.Range("F17").Value = 51
$ws.Range("G17").Value = 1419.75
$ws.Range("B18").Value = 667.5
$ws.Range("C18").Value = 647
$ws.Range("D18").Value = 647.70000000000005
$ws.Range("E18").Value = 649.1
$ws.Range("F18").Value = 24
$ws.Range("G18").Value = 664.3
$ws.Range("B19").Value = 426.1
$ws.Range("C19").Value = 419.9
$ws.Range("D19").Value = 420.45
$ws.Range("E19").Value = 421.55
$ws.Range("F19").Value = 15
$ws.Range("G19").Value = 425.95
$ws.Range("B20").Value = 1553.3
$ws.Range("C20").Value = 1528.05
$ws.Range("D20").Value = 1542
$ws.Range("E20").Value = 1542.2
$ws.Range("F20").Value = 21
$ws.Range("G20").Value = 1550.9
$ws.Range("B21").Value = 306
$ws.Range("C21").Value = 299
$ws.Range("D21").Value = 304.3
$ws.Range("E21").Value = 304.05
$ws.Range("F21").Value = 71
$ws.Range("G21").Value = 303
$ws.Range("B22").Value = 2539
$ws.Range("C22").Value = 2471
$ws.Range("D22").Value = 2475.8000000000002
$ws.Range("E22").Value = 2479.8000000000002
$ws.Range("F22").Value = 70
$ws.Range("G22").Value = 2529.4
$ws.Range("B23").Value = 582
$ws.Range("C23").Value = 575.54999999999995
$ws.Range("D23").Value = 577
$ws.Range("E23").Value = 576.95000000000005
$ws.Range("F23").Value = 241
$ws.Range("G23").Value = 580.35
$ws.Range("B24").Value = 585.9
$ws.Range("C24").Value = 566.04999999999995
$ws.Range("D24").Value = 582.70000000000005
$ws.Range("E24").Value = 582.54999999999995
$ws.Range("F24").Value = 46
$ws.Range("G24").Value = 571.1
$ws.Range("B25").Value = 1017.4
$ws.Range("C25").Value = 1004.15
$ws.Range("D25").Value = 1017
$ws.Range("E25").Value = 1015.65
$ws.Range("F25").Value = 5
$ws.Range("G25").Value = 1010
$ws.Range("B26").Value = 619
$ws.Range("C26").Value = 609.65
$ws.Range("D26").Value = 611.65
$ws.Range("E26").Value = 610.95000000000005
$ws.Range("F26").Value = 70
$ws.Range("G26").Value = 617.25
$ws.Range("B27").Value = 252.8
$ws.Range("C27").Value = 246.55
$ws.Range("D27").Value = 249.8
$ws.Range("E27").Value = 249.75
$ws.Range("F27").Value = 286
$ws.Range("G27").Value = 247.35
$ws.Range("B28").Value = 120.1
$ws.Range("C28").Value = 118.05
$ws.Range("D28").Value = 118.3
$ws.Range("E28").Value = 118.25
$ws.Range("F28").Value = 271
$ws.Range("G28").Value = 119.65
$ws.Range("B29").Value = 8324
$ws.Range("C29").Value = 8198
$ws.Range("D29").Value = 8198
$ws.Range("E29").Value = 8211.4
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 8213.2000000000007

# Update the active selection to match the new workbook state
$ws.Activate()
$ws.Range("J17").Select() | Out-Null
